# "Martin & Daniel nicknams & task."
# Adds a new "Villanova bug" task (priority A) to the Features sheet and
# assigns it to Martin & Daniel, plus bumps the indent on that sheet's
# task column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features")

# New row of data describing the task, its priority and who it's assigned to.
$ws.Range("A9").Value = "Villanova bug"
$ws.Range("B9").Value = "A"
$ws.Range("C9").Value = "Martin&Daniel"

# Increase the indent of the task column (was 7, now 9) to match the rest
# of the formatting used for this sheet.
$ws.Range("A2:A8").IndentLevel = 9

# Leave the cursor parked on the newly-entered cell, like Excel would after
# typing the data in by hand.
$null = $ws.Range("D9").Select()
